$d = $word.ActiveDocument

# 1. Update the "Nombre" cell text (CU-05 title) to reference CU-04
$d.Content.Find.Execute(
    "CU-05 Editar producto",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "CU-05 Editar producto extendido de CU-04 Ver producto",
    2
) | Out-Null

# 2. Add a new "Extiende" / "CU-04" row at the end of the table,
#    right after the "Postcondiciones" row.
$tbl = $d.Tables(1)
$newRow = $tbl.Rows.Add()
$newRow.Cells(1).Range.Text = "Extiende "
$newRow.Cells(1).Range.Font.Bold = 1
$newRow.Cells(2).Range.Text = "CU-04"
